# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.04763786555579896, 0.04240448674262143, 3.900430680208489, 8.660232485948974, 12.65070551845588)
    3  = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 5.553084769722144)
    4  = @(3.230985683306322, 114.8270160096505, 0.1575252929769615, 8.660232485948974, 126.8757594718828)
    5  = @(1.459612070389937, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 3.781711156805759)
    6  = @(1.459612070389937, 1.667794583268128, 0.8054896365839992, 8.660232485948974, 12.59312877619104)
    7  = @(0.6753301551942219, 3099.503889238888, 3.900430680208489, 8.660232485948974, 3112.73988256024)
    8  = @(0.04763786555579896, 0.3127903958511391, 3.900430680208489, 8.660232485948974, 12.9210914275644)
    9  = @(0.6753301551942219, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 3.645393585217082)
    10 = @(1.459612070389937, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 3.781711156805759)
    11 = @(1.459612070389937, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 4.429675500412797)
    12 = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 5.553084769722144)
    13 = @(1.459612070389937, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 4.429675500412797)
    14 = @(3.230985683306322, 1.667794583268128, 3.900430680208489, 0.496779210170732, 9.295990156953671)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 2).Value = $vals[0]
    $ws.Cells.Item($r, 3).Value = $vals[1]
    $ws.Cells.Item($r, 4).Value = $vals[2]
    $ws.Cells.Item($r, 5).Value = $vals[3]
    $ws.Cells.Item($r, 7).Value = $vals[4]
}
